$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.665.00"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").Value = "1.866.26"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.97"
$ws.Range("E5").Value = "  +3.20%  "

$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4681"
$ws.Range("E7").Value = "  +4.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  +2.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.26"
$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08034"
$ws.Range("E10").Value = "  +1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.021"
$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.76"
$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("D13").Value = "1.873.29"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.931"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.127"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.59"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06637"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").Value = "27.682.60"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").Value = "2.106.03"
$ws.Range("E26").Value = "  +1.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.41"
$ws.Range("E27").Value = "  +4.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.13"
$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.085"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.543"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.68"
$ws.Range("E31").Value = "  +2.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9651"
$ws.Range("E32").Value = "  +3.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09466"
$ws.Range("E33").Value = "  +2.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.442"
$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.597"
$ws.Range("E35").Value = "  +0.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.307"
$ws.Range("E36").Value = "  +0.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02250"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06071"
$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.232"
$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.105"
$ws.Range("E40").Value = "  -1.81%  "

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5967"
$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1893"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.21"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.267"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5674"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.24"
$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.387"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06845"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.02"
$ws.Range("E51").Value = "  +5.55%  "
